$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1, shifting existing data down
$ws.Rows("1:1").Insert()

# Set the new header value and make it bold
$ws.Range("A1").Value = "Sygedage "
$ws.Range("A1").Font.Bold = $true

# Update selection to match the target state
$ws.Range("E9").Select()
